# "Another try on vectors"
# Updates the Lawrence Lookups mean-vector table (M:N) with a new set of
# label/value pairs for rows 13-25, adds a trailing label in J43, and
# moves the active sheet/selection from "Lawernce Mean Vectors" to
# "Lawrence Lookups".

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Lawernce Mean Vectors")
$ws4 = $wb.Worksheets.Item("Lawrence Lookups")

# --- Update the M:N lookup pairs on the "Lawrence Lookups" sheet ---
# Each row keeps its N (score) value except where noted; the M (label)
# values are re-shuffled to new subjects.
$ws4.Range("M13").Value = "Political Science"
$ws4.Range("N13").Value = 1

$ws4.Range("M14").Value = "Literature"
$ws4.Range("N14").Value = 1

$ws4.Range("M15").Value = "Social Studies"
$ws4.Range("N15").Value = 0

$ws4.Range("M16").Value = "Psychology"
$ws4.Range("N16").Value = 0

$ws4.Range("M17").Value = "Sociology"
$ws4.Range("N17").Value = 0

$ws4.Range("M18").Value = "Anthropology"
$ws4.Range("N18").Value = 0

$ws4.Range("M19").Value = "Linguistics"
$ws4.Range("N19").Value = 0

$ws4.Range("M20").Value = "Management"
$ws4.Range("N20").Value = -1

$ws4.Range("M21").Value = "History"
$ws4.Range("N21").Value = -1

$ws4.Range("M22").Value = "Engineering"
$ws4.Range("N22").Value = -1

$ws4.Range("M23").Value = "Accounting"
$ws4.Range("N23").Value = -1

$ws4.Range("M24").Value = "Computer Science"
$ws4.Range("N24").Value = -2

$ws4.Range("M25").Value = "Math"
$ws4.Range("N25").Value = -2

# --- New trailing row ---
$ws4.Range("J43").Value = "s"

# --- Move the active sheet/selection from "Lawernce Mean Vectors" to
#     "Lawrence Lookups" ---
$ws3.Activate()
$ws3.Range("K33").Select()

$ws4.Activate()
$ws4.Range("J43").Select()
